# Update the "K" column (column G) with newly-regenerated strikeout values
# (commit message: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals"). Column G header is "K"; rows 2-36
# hold one game record each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(5,4,6,9,4,10,9,11,6,13,7,13,5,6,10,8,5,9,5,9,3,4,4,4,5,10,8,4,7,7,7,5,6,4,2)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
